# "Add files via upload" -- re-uploaded workbook fixes a data-entry slip on the
# five "actual vs predicted CO2" sheets (Cambodia, Laos, Myanmar, Thailand,
# Vietnam): for the forecast years 2015-2019 (rows 57-61) the predicted value
# had been typed into column B ("Actual CO2 emissions") instead of column C
# ("Predicted CO2 emissions"). Column B should read 0 (no actual data for a
# future year) and column C should carry the prediction.
#
# The author's last on-screen action also ends up with B57:C61 selected on
# each of those sheets (active cell B57), and Cambodia -- the first tab --
# left as the active sheet/tab when the file was saved.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Cambodia", "Laos", "Myanmar", "Thailand", "Vietnam")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    for ($r = 57; $r -le 61; $r++) {
        $predicted = $ws.Cells.Item($r, 2).Value2
        $ws.Cells.Item($r, 2).Value = 0
        $ws.Cells.Item($r, 3).Value = $predicted
    }
}

# Leave B57:C61 selected on every affected sheet. Visit Cambodia *last* so
# it ends up the active sheet/tab -- matching the saved file's
# activeTab/tabSelected state (activeTab moves from Vietnam (4) to
# Cambodia (0), i.e. the attribute is simply omitted since 0 is the
# schema default).
$selectOrder = @("Laos", "Myanmar", "Thailand", "Vietnam", "Cambodia")
foreach ($name in $selectOrder) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Select()
    $ws.Range("B57:C61").Select()
}

Write-Host "done"
